# "Stückliste_Drohne" BOM: fill in the real part quantities that were still
# marked with the "?" placeholder in the "Anzahl" column.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

$ws.Range("C15").Value = 2   # Reinigungsdüsen
$ws.Range("C26").Value = 1   # Antriebsrad
$ws.Range("C28").Value = 1   # Überträgt die Antriebsleistung...
$ws.Range("C30").Value = 1   # Laserscanner

# Move the selection/cursor to C13 (and drop the scrolled-away top-left cell).
$ws.Range("C13").Select()
